$d = $word.ActiveDocument

# Step 1: change the first run's text to add two trailing spaces
$d.Content.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false, $true, 1, $false, "This is a Microsoft word document.  ", 2)

# Step 2: append three new red-colored runs at the end of paragraph 1,
# right before the paragraph mark, forming:
#   "(This is a change " + [EN DASH] + " Version for main branch)"
$para1 = $d.Paragraphs(1)
$insertPos = $para1.Range.End - 1

$dash = [char]0x2013

$part1 = "(This is a change " + $dash + " Ve"
$part2 = "rsion for main branch"
$part3 = ")"

$r1 = $d.Range($insertPos, $insertPos)
$r1.InsertAfter($part1)
$r1.Font.Color = 255

$pos2 = $insertPos + $part1.Length
$r2 = $d.Range($pos2, $pos2)
$r2.InsertAfter($part2)
$r2.Font.Color = 255

$pos3 = $pos2 + $part2.Length
$r3 = $d.Range($pos3, $pos3)
$r3.InsertAfter($part3)
$r3.Font.Color = 255
